$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores all these cells as text (inline strings), including
# values that look numeric (e.g. "0.998", "172.27"). Force text format on the
# Price column before assignment so Excel does not silently convert them to
# numeric cells, which would change both the stored type and value precision.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.159.42"
$ws.Range("E2").Value = "  -3.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.745.96"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.43"
$ws.Range("E5").Value = "  -3.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.27"
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.748.08"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  -4.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.27"
$ws.Range("E11").Value = "  -3.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -3.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.67"
$ws.Range("E13").Value = "  -4.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  -4.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.370.34"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.729.74"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.254.76"
$ws.Range("E17").Value = "  -3.26%  "
$ws.Range("E18").Value = "  -4.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.08"
$ws.Range("E19").Value = "  -5.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.05"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "485.15"
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.07"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.720"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.74"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("E25").Value = "  -8.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000136"
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.18"
$ws.Range("E27").Value = "  -5.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  -9.37%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.90"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.40"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.06"
$ws.Range("E32").Value = "  +5.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.67"
$ws.Range("E33").Value = "  -3.75%  "
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -4.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.135"
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.69"
$ws.Range("E38").Value = "  -6.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.323"
$ws.Range("E39").Value = "  -6.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "449.35"
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "48.63"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.97"
$ws.Range("E42").Value = "  -3.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.85"
$ws.Range("E43").Value = "  -4.97%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.28"
$ws.Range("E44").Value = "  -8.68%  "
$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.21"
$ws.Range("E45").Value = "  -3.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.795.89"
$ws.Range("E46").Value = "  -5.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "139.48"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0347"
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.79"
$ws.Range("E50").Value = "  -4.18%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.09"
$ws.Range("E51").Value = "  +8.54%  "
